# RPA datasets push 2024-06-15
# Insert a new IPO record ("그리드위즈") at the top of the data table
# (row 2), shifting all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..14) down to (3..15), working bottom-up so
# we never overwrite a row before it has been copied.
# (NOTE: use Value2, not Value -- the latter's getter is unreliable here.)
for ($r = 15; $r -ge 3; $r--) {
    $src = $r - 1
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($src, $c).Value2
    }
}

# Populate the newly freed row 2 with the new subscription record.
$ws.Cells.Item(2, 1).Value2  = "2024-06-03"
$ws.Cells.Item(2, 2).Value2  = "그리드위즈"
$ws.Cells.Item(2, 3).Value2  = "삼성"
$ws.Cells.Item(2, 4).Value2  = "2024-06-07"
$ws.Cells.Item(2, 5).Value2  = "2024-06-14"
$ws.Cells.Item(2, 6).Value2  = 56000000
$ws.Cells.Item(2, 7).Value2  = 1400000
$ws.Cells.Item(2, 8).Value2  = "-"
$ws.Cells.Item(2, 9).Value2  = 34000
$ws.Cells.Item(2, 10).Value2 = 40000
$ws.Cells.Item(2, 11).Value2 = "-"
$ws.Cells.Item(2, 12).Value2 = 40000
$ws.Cells.Item(2, 13).Value2 = "-"
$ws.Cells.Item(2, 14).Value2 = "-"
$ws.Cells.Item(2, 15).Value2 = 0
$ws.Cells.Item(2, 16).Value2 = "-"
$ws.Cells.Item(2, 17).Value2 = "-"
$ws.Cells.Item(2, 18).Value2 = "569.89 : 1"
$ws.Cells.Item(2, 19).Value2 = "-"
$ws.Cells.Item(2, 20).Value2 = "-"
